$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.120.28"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.052.69"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.23"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.664"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.59"
$ws.Range("E7").Value = "  +6.14%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0783"
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.108"
$ws.Range("E11").Value = "  +1.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.96"
$ws.Range("E12").Value = "  +5.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.353.12"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.813"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("E15").Value = "  +6.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.056.67"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "37.129.07"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.81"
$ws.Range("E18").Value = "  +18.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.56"
$ws.Range("E19").Value = "  +3.11%  "
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.35"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "236.44"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("E25").Value = "  +11.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.20"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.23"
$ws.Range("E27").Value = "  +2.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.83"
$ws.Range("E28").Value = "  -0.95%  "
$ws.Range("E29").Value = "  +0.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.15"
$ws.Range("E30").Value = "  +9.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.72"
$ws.Range("E31").Value = "  +4.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0613"
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.49"
$ws.Range("E33").Value = "  +4.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0896"
$ws.Range("E34").Value = "  +4.19%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.23"
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("E37").Value = "  -2.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.108"
$ws.Range("E38").Value = "  +5.92%  "
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.17"
$ws.Range("E40").Value = "  +14.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.17"
$ws.Range("E41").Value = "  +29.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.59"
$ws.Range("E42").Value = "  -2.96%  "
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.13"
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "96.00"
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.44"
$ws.Range("E46").Value = "  +1.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.281.66"
$ws.Range("E47").Value = "  -0.95%  "
$ws.Range("E48").Value = "  -1.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.76"
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.240.31"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.47"
$ws.Range("E51").Value = "  -16.48%  "
